$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width: 31.3984375 -> 23.296875 (drop bestFit autosize) ---
$ws.Columns("B:B").ColumnWidth = 22.5

# --- Footer block restructure (rows 37-38) ---
# Old layout:
#   B37 = "Fuente: ARTF. Agencia Reguladora del Transporte Ferroviario."
#   E37 = "Ultima actualización: mayo 2024"           (style s=3)
#   E38 = "Dirección General de Planeación"            (style s=3)
# New layout:
#   B37 = "Actualización: mayo 2024."                  (style: left/top aligned, like B3)
#   B38 = "Fuente: ARTF. Agencia Reguladora del Transporte Ferroviario."
#   E38 = empty cell, keeps style s=3

# Capture the format used on B3 (fontId=1 / vertical=top) and apply it to B37,
# then add horizontal=left so it matches the new style exactly.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B37").PasteSpecial(-4122) | Out-Null
$ws.Range("B37").HorizontalAlignment = -4131

# New text content for B37 / B38
$ws.Range("B37").Value = "Actualización: mayo 2024."
$ws.Range("B38").Value = "Fuente: ARTF. Agencia Reguladora del Transporte Ferroviario."

# E37 is fully removed (no value, no format) while E38 keeps its format
# but loses its value.
$ws.Range("E37").Clear() | Out-Null
$ws.Range("E38").ClearContents() | Out-Null

$excel.CutCopyMode = 0

# --- Title (fix accent: "ferreas" -> "férreas") ---
$ws.Range("B2").Value = "Longitud de vías férreas 2023"
